$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new sentence " Além disso, rotinas específica em linguagem C
#    foram desenvolvidas " right after "...chapas piezoelétricas." and before
#    "Este código foi desenvolvido...".
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    ". Este código foi desenvolvido seguindo os padrões SOLID e API ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Além disso, rotinas específica em linguagem C foram desenvolvidas  Este código foi desenvolvido seguindo os padrões SOLID e API ",
    2)
Write-Output ("inserted sentence: " + $found)

# ---------------------------------------------------------------------------
# 2. Relocate the existing "_GoBack" bookmark from the end of the document
#    (right after the lone page break) to right after the word we just
#    inserted ("desenvolvidas "), i.e. right before " Este código...".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$bmRng = $d.Content
$bmFound = $bmRng.Find.Execute("desenvolvidas ", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
Write-Output ("located insertion point for bookmark: " + $bmFound)
$bmRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Output "done"
